$p = $ppt.ActivePresentation

# Slide 10 ("Try test with shell") is the slide touched by this commit
# (cId="32353519" sldId="291" in the changesInfo part).
$s = $p.Slides.Item(10)

# --- Content Placeholder 2 (shape id=3): reword the bullet and give it an
#     explicit position/size (it previously inherited the layout's geometry).
$body = $s.Shapes.Item(2)
$body.Left = 84.24
$body.Top = 167.04
$body.Width = 792
$body.Height = 54.96
# Reset to a placeholder value first so PowerPoint collapses the text into a
# single run instead of reusing the old run boundaries.
$body.TextFrame.TextRange.Text = "x"
$body.TextFrame.TextRange.Text = "Use shell to start python interactive environment, and run some test code"

# --- Content Placeholder 2 (shape id=4, the ">>>" python shell textbox):
#     just reposition it.
$pyShell = $s.Shapes.Item(3)
$pyShell.Left = 203.25
$pyShell.Top = 255.8228

# --- Content Placeholder 2 (shape id=5, the "runserver" textbox): removed
#     entirely now that the slide only talks about one shell.
$s.Shapes.Item(4).Delete()
